# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from 2023-10-22 (45221) to 2023-10-25 (45224).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C15").Value = 45224
